$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Rows(10).Delete()
$ws1.Range("B3").Value = "2024-08-17"
$ws1.Range("C3").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws1.Range("D3").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E3").Value = "2024.08.17 09:30-08.18 17:00"
$ws1.Range("F3").Value = 5148
$ws1.Range("G3").Value = 69
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

$ws1.Range("B4").Value = "2024-08-17"
$ws1.Range("C4").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws1.Range("D4").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E4").Value = "2024.08.17 09:30-08.17 17:00"
$ws1.Range("F4").Value = 377
$ws1.Range("G4").Value = "已售罄"
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

$ws1.Range("B5").Value = "2024-08-17"
$ws1.Range("C5").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws1.Range("D5").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E5").Value = "2024.08.17 09:30-08.17 17:00"
$ws1.Range("F5").Value = 630
$ws1.Range("G5").Value = "已售罄"
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws1.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

$ws1.Range("B6").Value = "2024-08-17"
$ws1.Range("C6").Value = "合肥·银魂主题派对only2.0"
$ws1.Range("D6").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws1.Range("E6").Value = "2024.08.17 13:00-08.17 18:00"
$ws1.Range("F6").Value = 300
$ws1.Range("G6").Value = 128
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws1.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

$ws1.Range("B7").Value = "2024-08-18"
$ws1.Range("C7").Value = "合肥·SSS第五人格only"
$ws1.Range("D7").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws1.Range("E7").Value = "2024.08.18 09:00-08.18 17:00"
$ws1.Range("F7").Value = 779
$ws1.Range("G7").Value = 68
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws1.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"

$ws1.Range("B8").Value = "2024-09-07"
$ws1.Range("C8").Value = "合肥·国乙only宇宙心动（含夜场）"
$ws1.Range("D8").Value = "文忠路1865号 赫拉诺言艺术中心"
$ws1.Range("E8").Value = "2024.09.07 10:00-09.07 21:00"
$ws1.Range("F8").Value = 261
$ws1.Range("G8").Value = 48
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=89803"
$ws1.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"

$ws1.Range("B9").Value = "2024-09-16"
$ws1.Range("C9").Value = "肥西·星域动漫游戏嘉年华"
$ws1.Range("D9").Value = "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)"
$ws1.Range("E9").Value = "2024.09.16 10:00-09.16 17:00"
$ws1.Range("F9").Value = 7
$ws1.Range("G9").Value = 45
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=90489"
$ws1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows(12).Delete()
$ws4.Range("B3").Value = "2024-08-17"
$ws4.Range("C3").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws4.Range("D3").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E3").Value = "2024.08.17 09:30-08.18 17:00"
$ws4.Range("F3").Value = 5148
$ws4.Range("G3").Value = 69
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws4.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

$ws4.Range("B4").Value = "2024-08-17"
$ws4.Range("C4").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws4.Range("D4").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E4").Value = "2024.08.17 09:30-08.17 17:00"
$ws4.Range("F4").Value = 377
$ws4.Range("G4").Value = "已售罄"
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws4.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

$ws4.Range("B5").Value = "2024-08-17"
$ws4.Range("C5").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws4.Range("D5").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E5").Value = "2024.08.17 09:30-08.17 17:00"
$ws4.Range("F5").Value = 630
$ws4.Range("G5").Value = "已售罄"
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws4.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

$ws4.Range("B6").Value = "2024-08-17"
$ws4.Range("C6").Value = "合肥·银魂主题派对only2.0"
$ws4.Range("D6").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws4.Range("E6").Value = "2024.08.17 13:00-08.17 18:00"
$ws4.Range("F6").Value = 300
$ws4.Range("G6").Value = 128
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

$ws4.Range("B7").Value = "2024-08-18"
$ws4.Range("C7").Value = "合肥·SSS第五人格only"
$ws4.Range("D7").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws4.Range("E7").Value = "2024.08.18 09:00-08.18 17:00"
$ws4.Range("F7").Value = 779
$ws4.Range("G7").Value = 68
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"

$ws4.Range("B8").Value = "2024-08-25"
$ws4.Range("C8").Value = "合肥·CrossingX意次元｜乐队番ONLY同人"
$ws4.Range("D8").Value = "国祯广场B-1楼 背影骑士LIVEHOUSE"
$ws4.Range("E8").Value = "2024.08.25 13:30-08.25 16:00"
$ws4.Range("F8").Value = 32
$ws4.Range("G8").Value = 38
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=90032"
$ws4.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202407/GYPAnumr1721896597593.jpeg"

$ws4.Range("B9").Value = "2024-09-07"
$ws4.Range("C9").Value = "合肥·国乙only宇宙心动（含夜场）"
$ws4.Range("D9").Value = "文忠路1865号 赫拉诺言艺术中心"
$ws4.Range("E9").Value = "2024.09.07 10:00-09.07 21:00"
$ws4.Range("F9").Value = 256
$ws4.Range("G9").Value = 48
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=89803"
$ws4.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"

$ws4.Range("B10").Value = "2024-09-16"
$ws4.Range("C10").Value = "肥西·星域动漫游戏嘉年华"
$ws4.Range("D10").Value = "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)"
$ws4.Range("E10").Value = "2024.09.16 10:00-09.16 17:00"
$ws4.Range("F10").Value = 7
$ws4.Range("G10").Value = 45
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=90489"
$ws4.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"

$ws4.Range("B11").Value = "2024-10-26"
$ws4.Range("C11").Value = "合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$ws4.Range("D11").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws4.Range("E11").Value = "2024.10.26 19:30-10.26 21:00"
$ws4.Range("F11").Value = 5
$ws4.Range("G11").Value = 40
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90322"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg"

Write-Output "done"